$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Group A: Delhi_NewZealand) - fix duplicate match counts
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2

# Row 3 (Group B)
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 4

# Row 4 (Group B: RoyalChallengers_Bangladesh)
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 4

# Row 5 (Group A: Rajastan_Australia)
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4
